# Refresh the cryptos worksheet: latest Price (D) and Volume(1h) (E) values
# for each coin row, as pulled by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.848.58"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.316.01"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'97.38"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "'272.95"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "'45.44"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'0.0954"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "'7.98"
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "2.654.17"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "'15.52"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "'0.874"
$ws.Range("E16").Value = "  +7.41%  "
$ws.Range("D17").Value = "2.321.08"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").Value = "43.784.01"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("D20").Value = "'6.43"
$ws.Range("E20").Value = "  +4.61%  "
$ws.Range("D21").Value = "'73.56"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").Value = "'239.83"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'2.55"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "'11.38"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "'38.30"
$ws.Range("E30").Value = "  -6.25%  "
$ws.Range("D31").Value = "'22.46"
$ws.Range("E31").Value = "  +6.82%  "
$ws.Range("D32").Value = "'175.22"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "'5.50"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("D36").Value = "'0.0365"
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").Value = "'0.110"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "'4.47"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  -6.60%  "
$ws.Range("D40").Value = "'0.245"
$ws.Range("E40").Value = "  +8.16%  "
$ws.Range("D41").Value = "'2.41"
$ws.Range("E41").Value = "  +11.10%  "
$ws.Range("E42").Value = "  +23.19%  "
$ws.Range("D43").Value = "'12.40"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").Value = "'62.95"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "'9.22"
$ws.Range("E45").Value = "  +10.01%  "
$ws.Range("D46").Value = "'5.35"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "'1.21"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +15.59%  "
$ws.Range("D51").Value = "2.541.64"
$ws.Range("E51").Value = "  +3.29%  "
